$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 16453.5
$ws.Range("I20").Value = 15232.571
$ws.Range("J20").Value = 25000
$ws.Range("K20").Value = 15232.571
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = -15002.571
$ws.Range("N20").Value = -25460

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 16453.5
$ws.Range("I35").Value = 15232.571
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 15232.571
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -14853.571
$ws.Range("N35").Value = -25758

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 693.4545000000001
$ws.Range("I41").Value = 592.9
$ws.Range("J41").Value = 1699
$ws.Range("K41").Value = 592.9
$ws.Range("L41").Value = 1699
$ws.Range("M41").Value = -152.9
$ws.Range("N41").Value = -2579

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1210.25
$ws.Range("I96").Value = 2477.6667
$ws.Range("J96").Value = 449.8
$ws.Range("K96").Value = 7433.000100000001
$ws.Range("L96").Value = 1349.4
$ws.Range("M96").Value = -6060.000100000001
$ws.Range("N96").Value = -4095.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5021.1665
$ws.Range("I100").Value = 3539.6
$ws.Range("J100").Value = 6079.4287
$ws.Range("K100").Value = 3539.6
$ws.Range("L100").Value = 6079.4287
$ws.Range("M100").Value = -2998.6
$ws.Range("N100").Value = -7161.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2939.375
$ws.Range("I106").Value = 1757.75
$ws.Range("K106").Value = 1757.75
$ws.Range("M106").Value = -1126.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2927.7974
$ws.Range("I137").Value = 2155.4
$ws.Range("J137").Value = 3189.6272
$ws.Range("K137").Value = 6466.200000000001
$ws.Range("L137").Value = 9568.881600000001
$ws.Range("M137").Value = -3916.200000000001
$ws.Range("N137").Value = -14668.8816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5413.0415
$ws.Range("J2").Value = 10038.454
$ws.Range("L2").Value = 10038.454
$ws.Range("N2").Value = -10264.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 550.25
$ws.Range("I4").Value = 550.25
$ws.Range("K4").Value = 550.25
$ws.Range("M4").Value = -434.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4214.7163
$ws.Range("I32").Value = 4214.7163
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4214.7163
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3927.7163

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3890.5881
$ws.Range("I45").Value = 3535.75
$ws.Range("J45").Value = 4742.2
$ws.Range("K45").Value = 3535.75
$ws.Range("L45").Value = 4742.2
$ws.Range("M45").Value = -3158.75
$ws.Range("N45").Value = -5496.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2250.074
$ws.Range("I74").Value = 1623.3914
$ws.Range("K74").Value = 1623.3914
$ws.Range("M74").Value = -749.3914

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2250.074
$ws.Range("I77").Value = 1623.3914
$ws.Range("K77").Value = 8116.957
$ws.Range("M77").Value = -3748.957

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1053.0322
$ws.Range("I97").Value = 894.2
$ws.Range("K97").Value = 894.2
$ws.Range("M97").Value = -398.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 78125
$ws.Range("J107").Value = 78125
$ws.Range("L107").Value = 78125
$ws.Range("N107").Value = -85805

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 42795.668
$ws.Range("J112").Value = 42795.668
$ws.Range("L112").Value = 42795.668
$ws.Range("N112").Value = -45749.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5413.0415
$ws.Range("J116").Value = 10038.454
$ws.Range("L116").Value = 10038.454
$ws.Range("N116").Value = -14626.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4017.3333
$ws.Range("I122").Value = 3679.9285
$ws.Range("K122").Value = 11039.7855
$ws.Range("M122").Value = -8589.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3346.1738
$ws.Range("I132").Value = 2616.5715
$ws.Range("K132").Value = 7849.7145
$ws.Range("M132").Value = -5319.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 53824.617
$ws.Range("J135").Value = 53824.617
$ws.Range("L135").Value = 53824.617
$ws.Range("N135").Value = -63964.617

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 82429
$ws.Range("J138").Value = 82429
$ws.Range("L138").Value = 82429
$ws.Range("N138").Value = -92709

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5413.0415
$ws.Range("J3").Value = 10038.454
$ws.Range("L3").Value = 10038.454
$ws.Range("N3").Value = -10266.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4299.44
$ws.Range("I86").Value = 3399.158
$ws.Range("J86").Value = 7150.3335
$ws.Range("K86").Value = 3399.158
$ws.Range("L86").Value = 7150.3335
$ws.Range("M86").Value = -2276.158
$ws.Range("N86").Value = -9396.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4299.44
$ws.Range("I89").Value = 3399.158
$ws.Range("J89").Value = 7150.3335
$ws.Range("K89").Value = 16995.79
$ws.Range("L89").Value = 35751.6675
$ws.Range("M89").Value = -11379.79
$ws.Range("N89").Value = -46983.6675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2713.7144
$ws.Range("I99").Value = 1874.375
$ws.Range("J99").Value = 3832.8333
$ws.Range("K99").Value = 1874.375
$ws.Range("L99").Value = 3832.8333
$ws.Range("M99").Value = -376.375
$ws.Range("N99").Value = -6828.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 49613.332
$ws.Range("J112").Value = 49613.332
$ws.Range("L112").Value = 49613.332
$ws.Range("N112").Value = -52567.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3980.7144
$ws.Range("I134").Value = 3739.182
$ws.Range("K134").Value = 11217.546
$ws.Range("M134").Value = -8682.545999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34596.547
$ws.Range("I31").Value = 2090.818
$ws.Range("K31").Value = 2090.818
$ws.Range("M31").Value = -1795.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 34596.547
$ws.Range("I34").Value = 2090.818
$ws.Range("K34").Value = 2090.818
$ws.Range("M34").Value = -1888.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4391.2856
$ws.Range("I58").Value = 2683.8696
$ws.Range("K58").Value = 2683.8696
$ws.Range("M58").Value = -2480.8696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4115.4546
$ws.Range("I99").Value = 2999
$ws.Range("K99").Value = 2999
$ws.Range("M99").Value = -1501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3016.3333
$ws.Range("I107").Value = 2355.875
$ws.Range("J107").Value = 3771.1428
$ws.Range("K107").Value = 2355.875
$ws.Range("L107").Value = 3771.1428
$ws.Range("M107").Value = -435.875
$ws.Range("N107").Value = -7611.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4115.4546
$ws.Range("I126").Value = 2999
$ws.Range("K126").Value = 8997
$ws.Range("M126").Value = -6527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5056.5557
$ws.Range("I134").Value = 3749.1667
$ws.Range("J134").Value = 7671.3335
$ws.Range("K134").Value = 11247.5001
$ws.Range("L134").Value = 23014.0005
$ws.Range("M134").Value = -8712.500100000001
$ws.Range("N134").Value = -28084.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4391.2856
$ws.Range("I136").Value = 2683.8696
$ws.Range("K136").Value = 8051.6088
$ws.Range("M136").Value = -5501.6088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 217842.17
$ws.Range("J141").Value = 227246.4
$ws.Range("L141").Value = 227246.4
$ws.Range("N141").Value = -237606.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1831.1578
$ws.Range("I107").Value = 1692.6428
$ws.Range("J107").Value = 2219
$ws.Range("K107").Value = 5077.928400000001
$ws.Range("L107").Value = 6657
$ws.Range("M107").Value = -3157.928400000001
$ws.Range("N107").Value = -10497

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3093.12
$ws.Range("I132").Value = 2427.087
$ws.Range("J132").Value = 10752.5
$ws.Range("K132").Value = 7281.261
$ws.Range("L132").Value = 32257.5
$ws.Range("M132").Value = -4751.261
$ws.Range("N132").Value = -37317.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2741.8438
$ws.Range("I136").Value = 1197.4073
$ws.Range("J136").Value = 11081.8
$ws.Range("K136").Value = 3592.2219
$ws.Range("L136").Value = 33245.39999999999
$ws.Range("M136").Value = -1042.2219
$ws.Range("N136").Value = -38345.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 82429
$ws.Range("J140").Value = 82429
$ws.Range("L140").Value = 82429
$ws.Range("N140").Value = -92789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 104647.336
$ws.Range("J141").Value = 104647.336
$ws.Range("L141").Value = 104647.336
$ws.Range("N141").Value = -115007.336
